$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G15").Value = '2025/2026, 2023/2024'
$ws.Range("G19").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G21").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G22").Value = '2025/2026, Eman_mohamed@med.asu.edu.eg'
$ws.Range("G37").Value = '2025/2026, 2023/2024'
$ws.Range("G41").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G43").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G44").Value = '2025/2026, Eman_mohamed@med.asu.edu.eg'
$ws.Range("G60").Value = '2025/2026, 2026/2027'
$ws.Range("G63").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G64").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G65").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G66").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G82").Value = '2025/2026, 2026/2027'
$ws.Range("G85").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G86").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G87").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G88").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G89").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G104").Value = '2025/2026, 2026/2027'
$ws.Range("G106").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G107").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G110").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G111").Value = '2025/2026, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G126").Value = '2025/2026, 2026/2027'
$ws.Range("G128").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G129").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G132").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G133").Value = '2025/2026, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G150").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G153").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G155").Value = '2025/2026, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G172").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G175").Value = '2025/2026, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G177").Value = '2025/2026, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'